$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (post-edit) data for rows 2..19, columns A..E.
# A = date_of_forecast (serial date, style carried from existing date cells)
# B = y_0 (year)
# C = y_0_forecast
# D = y_1 (year)
# E = y_1_forecast (blank for the first few rows)
$data = @(
    @(39400, 2007, 1.75539628881467,  2008, $null),
    @(39765, 2008, 2.213911448916162, 2009, $null),
    @(40130, 2009, 2.533533936850563, 2010, $null),
    @(40494, 2010, 2.088987486264915, 2011, 1.485473821631844),
    @(40862, 2011, 1.212544822741002, 2012, 1.799394172339341),
    @(41228, 2012, 1.196776590518644, 2013, 1.2151583353186),
    @(41592, 2013, 0.4712609263772594, 2014, 1.107727073902187),
    @(41957, 2014, 0.8783377572271434, 2015, 1.612081704302182),
    @(42321, 2015, 2.29066283401107,  2016, 2.221748592150097),
    @(42689, 2016, 4.109890522944348, 2017, 2.932944072183674),
    @(43053, 2017, 1.336316831462692, 2018, 1.104283769064729),
    @(43418, 2018, 1.197912858979611, 2019, 1.649865498505276),
    @(43783, 2019, 1.727537197898665, 2020, 2.284828905445169),
    @(44159, 2020, 3.647228437274408, 2021, 3.474365686630398),
    @(44525, 2021, 2.777797690741424, 2022, 1.742844348069261),
    @(44890, 2022, 0.6994919452575576, 2023, 0.5651273241891186),
    @(45254, 2023, -1.432689847121871, 2024, 0.4518870186319468),
    @(45618, 2024, 2.033479419175133, 2025, 1.959987726090251)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    if ($row[4] -ne $null) {
        $ws.Cells.Item($r, 5).Value2 = $row[4]
    }
    $r = $r + 1
}

# Make sure every date cell in column A (rows 2..19) carries the same
# style as the original date column (copy format from a row that already
# had it, so no new style entries are introduced).
$ws.Range("A2").Copy()
$ws.Range("A3:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
